$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at row 1127, shifting existing rows 1127:1204 down to 1131:1208
$ws.Range("A1127:A1130").EntireRow.Insert()

# Row 1127
$ws.Cells.Item(1127, "A").Value = 3
$ws.Cells.Item(1127, "B").Value = 'Femacal de La Calera'
$ws.Cells.Item(1127, "C").Value = 'Coquimbo'
$ws.Cells.Item(1127, "D").Value = 44585
$ws.Cells.Item(1127, "E").Value = 5
$ws.Cells.Item(1127, "F").Value = 100112002
$ws.Cells.Item(1127, "G").Value = 'Pimiento'
$ws.Cells.Item(1127, "H").Value = 'Zafiro amarillo'
$ws.Cells.Item(1127, "I").Value = 'Primera'
$ws.Cells.Item(1127, "J").Value = 73
$ws.Cells.Item(1127, "K").Value = 25000
$ws.Cells.Item(1127, "L").Value = 26000
$ws.Cells.Item(1127, "M").Value = 25479
$ws.Cells.Item(1127, "N").Value = '$/caja 18 kilos'
$ws.Cells.Item(1127, "O").Value = 'Limache'
$ws.Cells.Item(1127, "P").Value = 1416
$ws.Cells.Item(1127, "Q").Value = 18
$ws.Cells.Item(1127, "R").Value = 'Hortaliza'

# Row 1128
$ws.Cells.Item(1128, "A").Value = 3
$ws.Cells.Item(1128, "B").Value = 'Femacal de La Calera'
$ws.Cells.Item(1128, "C").Value = 'Coquimbo'
$ws.Cells.Item(1128, "D").Value = 44585
$ws.Cells.Item(1128, "E").Value = 5
$ws.Cells.Item(1128, "F").Value = 100112002
$ws.Cells.Item(1128, "G").Value = 'Pimiento'
$ws.Cells.Item(1128, "H").Value = 'Zafiro rojo'
$ws.Cells.Item(1128, "I").Value = 'Primera'
$ws.Cells.Item(1128, "J").Value = 73
$ws.Cells.Item(1128, "K").Value = 23000
$ws.Cells.Item(1128, "L").Value = 24000
$ws.Cells.Item(1128, "M").Value = 23479
$ws.Cells.Item(1128, "N").Value = '$/caja 18 kilos'
$ws.Cells.Item(1128, "O").Value = 'Limache'
$ws.Cells.Item(1128, "P").Value = 1304
$ws.Cells.Item(1128, "Q").Value = 18
$ws.Cells.Item(1128, "R").Value = 'Hortaliza'

# Row 1129
$ws.Cells.Item(1129, "A").Value = 3
$ws.Cells.Item(1129, "B").Value = 'Femacal de La Calera'
$ws.Cells.Item(1129, "C").Value = 'Coquimbo'
$ws.Cells.Item(1129, "D").Value = 44585
$ws.Cells.Item(1129, "E").Value = 5
$ws.Cells.Item(1129, "F").Value = 100112002
$ws.Cells.Item(1129, "G").Value = 'Pimiento'
$ws.Cells.Item(1129, "H").Value = 'Zafiro verde'
$ws.Cells.Item(1129, "I").Value = 'Primera'
$ws.Cells.Item(1129, "J").Value = 73
$ws.Cells.Item(1129, "K").Value = 16000
$ws.Cells.Item(1129, "L").Value = 17000
$ws.Cells.Item(1129, "M").Value = 16521
$ws.Cells.Item(1129, "N").Value = '$/caja 18 kilos'
$ws.Cells.Item(1129, "O").Value = 'Limache'
$ws.Cells.Item(1129, "P").Value = 918
$ws.Cells.Item(1129, "Q").Value = 18
$ws.Cells.Item(1129, "R").Value = 'Hortaliza'

# Row 1130
$ws.Cells.Item(1130, "A").Value = 3
$ws.Cells.Item(1130, "B").Value = 'Femacal de La Calera'
$ws.Cells.Item(1130, "C").Value = 'Coquimbo'
$ws.Cells.Item(1130, "D").Value = 44585
$ws.Cells.Item(1130, "E").Value = 5
$ws.Cells.Item(1130, "F").Value = 100112002
$ws.Cells.Item(1130, "G").Value = 'Pimiento'
$ws.Cells.Item(1130, "H").Value = 'Zafiro verde'
$ws.Cells.Item(1130, "I").Value = 'Segunda'
$ws.Cells.Item(1130, "J").Value = 38
$ws.Cells.Item(1130, "K").Value = 13000
$ws.Cells.Item(1130, "L").Value = 13000
$ws.Cells.Item(1130, "M").Value = 13000
$ws.Cells.Item(1130, "N").Value = '$/caja 18 kilos'
$ws.Cells.Item(1130, "O").Value = 'Limache'
$ws.Cells.Item(1130, "P").Value = 722
$ws.Cells.Item(1130, "Q").Value = 18
$ws.Cells.Item(1130, "R").Value = 'Hortaliza'
